$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "IP00034"
$ws.Range("A3").Value = "IP00034"
$ws.Range("A1").Value = "code_import"
$ws.Range("B1").Value = "id_material_detail"

$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()

$ws.Range("E9").Select()
